$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Vegetation Tyoe Grassland" sheet right before
#    "Perennial Cropland Tyoe" (i.e. right after "Age Classes").
# ------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("Perennial Cropland Tyoe")
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "Vegetation Tyoe Grassland"

# ------------------------------------------------------------------
# 2. Populate the data.
# ------------------------------------------------------------------
$headers = @("ID", "vegetation_type", "description", "ratio_bgb_agb", "n", "range")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @(1, "Steppe",       "ST", 4,   7,  "±150%"),
    @(2, "Tundra",       "TU", 4,   7,  "±150%"),
    @(3, "Prairie",      "PR", 4,   7,  "±150%"),
    @(4, "Semi-Arid",    "SA", 2.8, 9,  "±95%"),
    @(5, "Sub-Tropical", "ST", 1.6, 7,  "±130%"),
    @(6, "Tropical",     "TR", 1.6, 7,  "±130%"),
    @(7, "Woodland",     "WL", 0.5, 19, "±80%"),
    @(8, "Savannah",     "SV", 0.5, 19, "±80%"),
    @(9, "Shrubland",    "SH", 2.8, 9,  "±144%")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ------------------------------------------------------------------
# 3. Formatting: reuse the header/body styles already present on the
#    "Perennial Cropland Tyoe" sheet (blue header band + thin borders)
#    so the workbook does not grow new, duplicate style records.
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("Perennial Cropland Tyoe")

# Header row: columns A & D use the "outer" style, B & C the "inner" one.
$srcSheet.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$srcSheet.Range("B1:C1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Body rows: columns A & D use the "outer" style, B, C, E & F the "inner" one.
$srcSheet.Range("A2:D2").Copy()
$ws.Range("A2:D10").PasteSpecial(-4122)
$srcSheet.Range("B2:C2").Copy()
$ws.Range("E2:F10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4. View state tweaks captured in the diff.
# ------------------------------------------------------------------
$ws.Range("A2:D10").Select()

$wsGrowing = $wb.Worksheets.Item("Growing stock level")
$wsGrowing.Range("C19").Select()

$wsAge = $wb.Worksheets.Item("Age Classes")
$wsAge.Range("B2").Select()

$wsContinent = $wb.Worksheets.Item("Continent Type")
$excel.ActiveWindow.ScrollColumn = 4
$wsContinent.Range("S41").Select()

$wsLandSub = $wb.Worksheets.Item("LandSubcategories")
$wsLandSub.Select()
